$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")
$ws.Activate()

# --- Rename "Team Grubi" -> "Team Grabi" for the two header-group rows that
# use it on this sheet (row 4 and row 10, both use style s=11) ---
$ws.Range("A4").Value = "Team Grabi"
$ws.Range("A10").Value = "Team Grabi"

# --- Update the competition name (column B) for the two existing blocks ---
$ws.Range("B2:B7").Value = "Fußballgolf"
$ws.Range("B8:B13").Value = "PAM"

# --- Build the three new 6-row blocks (rows 14-31) by duplicating the
# existing formatted block (A2:C7) so the same header/body cell styles
# (s=10 on the first row, s=11 on the rest) carry over, then overwrite the
# competition name + score values for each new block. ---
$ws.Range("A2:C7").Copy($ws.Range("A14:C19"))
$ws.Range("A2:C7").Copy($ws.Range("A20:C25"))
$ws.Range("A2:C7").Copy($ws.Range("A26:C31"))

# Fix the third row of each copied block back to "Team Grabi" (copy above
# reproduced the current A4 text, but make sure it is explicit/correct).
$ws.Range("A16").Value = "Team Grabi"
$ws.Range("A22").Value = "Team Grabi"
$ws.Range("A28").Value = "Team Grabi"

# Block: Tic-Tac-Toe (rows 14-19)
$ws.Range("B14:B19").Value = "Tic-Tac-Toe"
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 10
$ws.Range("C19").Value = 8

# NOTE: "Backbacking" (row block 26-31) is entered before "Ball hochhalten"
# (row block 20-25) so the shared-string table picks up the same ordering
# (57="Backbacking", 58="Ball hochhalten") the author's workbook has, even
# though that block's rows sit further down the sheet.

# Block: Backbacking (rows 26-31)
$ws.Range("B26:B31").Value = "Backbacking"
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 6
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 4
$ws.Range("C30").Value = 8
$ws.Range("C31").Value = 2

# --- Best-fit column B now, while the longest entry is still one of the
# 11-character competition names ("Fußballgolf"/"Tic-Tac-Toe"/"Backbacking"),
# matching the width the author's workbook ended up with. ---
$ws.Columns.Item(2).AutoFit()

# Block: Ball hochhalten (rows 20-25)
$ws.Range("B20:B25").Value = "Ball hochhalten"
$ws.Range("C20").Value = 10
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = 8
$ws.Range("C24").Value = 2
$ws.Range("C25").Value = 6

# --- Match the author's final cursor position/selection ---
$ws.Range("F16").Select()
